$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert existing Company Number cells (B2:B10) from text to real numbers
$ws.Range("B2").Value = 16455471
$ws.Range("B3").Value = 16455405
$ws.Range("B4").Value = 16455468
$ws.Range("B5").Value = 16455594
$ws.Range("B6").Value = 16455597
$ws.Range("B7").Value = 16455494
$ws.Range("B8").Value = 16455573
$ws.Range("B9").Value = 16455443
$ws.Range("B10").Value = 16455528

# New rows appended (11-23), Company Number (and date-like columns) stored as
# plain text, matching how the source data originally looked.
$newRows = @(
    @("LENDING CONSULTANCY LTD", "16455471", "2025-05-18", "active", "SIC", "2025-05-18", "2025-05-18 21:42:48"),
    @("ECHO VENTURES GROUP LIMITED", "16455744", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:48"),
    @("ESLB INVESTMENTS LIMITED", "16455669", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:48"),
    @("JISA VENTURES LTD", "16455405", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:49"),
    @("TALLY M E VENTURES LIMITED", "16455468", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:49"),
    @("PERFICIENT VENTURES LTD", "16455594", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:49"),
    @("BLUEBOW TECHNOLOGIES LTD", "16455597", "2025-05-18", "active", "SIC", "2025-05-18", "2025-05-18 21:42:49"),
    @("JENKINS VENTURES LTD", "16455788", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:49"),
    @("BIEN DEVELOPMENTS LTD", "16455494", "2025-05-18", "active", "SIC", "2025-05-18", "2025-05-18 21:42:50"),
    @("ALPHA HAULAGE SOLUTIONS LTD", "16455573", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:50"),
    @("MARIOS PROPERTY INVESTMENTS LTD", "16455816", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:50"),
    @("MARKOVIAN INVESTMENTS LIMITED", "16455443", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:50"),
    @("PARTNERS AMERICAN WHISKEY LTD", "16455528", "2025-05-18", "active", "Keyword", "2025-05-18", "2025-05-18 21:42:50")
)

$startRow = 11
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Force text format on columns that would otherwise be auto-detected by
    # Excel as numbers (B) or dates (C, F) so the values stay plain text.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 6).NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
